$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.713252999999999
$ws.Cells.Item(2, 8).Value = 20.139759
$ws.Cells.Item(2, 9).Value = 0.3101840064655811
$ws.Cells.Item(2, 10).Value = 0.3231642354899327
$ws.Cells.Item(2, 13).Value = 44.04223000000001
$ws.Cells.Item(2, 14).Value = 132.12669
$ws.Cells.Item(2, 15).Value = 0.1792438957545786
$ws.Cells.Item(2, 16).Value = 0.1868246872369915
$ws.Cells.Item(2, 17).Value = 295.6666326741901
$ws.Cells.Item(2, 18).Value = 2660.99969406771
$ws.Cells.Item(2, 19).Value = 0.05559858971965417
$ws.Cells.Item(2, 20).Value = 0.06037505722158815
$ws.Cells.Item(3, 7).Value = 6.713252999999999
$ws.Cells.Item(3, 8).Value = 20.139759
$ws.Cells.Item(3, 9).Value = 0.3101840064655811
$ws.Cells.Item(3, 10).Value = 0.3231642354899327
$ws.Cells.Item(3, 15).Value = 0.2018201397722426
$ws.Cells.Item(3, 16).Value = 0.2103557520458098
$ws.Cells.Item(3, 17).Value = 332.9066291551469
$ws.Cells.Item(3, 18).Value = 2996.159662396322
$ws.Cells.Item(3, 19).Value = 0.06260137953999778
$ws.Cells.Item(3, 20).Value = 0.06797945579079397
$ws.Cells.Item(4, 7).Value = 6.713252999999999
$ws.Cells.Item(4, 8).Value = 20.139759
$ws.Cells.Item(4, 9).Value = 0.3101840064655811
$ws.Cells.Item(4, 10).Value = 0.3231642354899327
$ws.Cells.Item(4, 13).Value = 51.56497066666667
$ws.Cells.Item(4, 14).Value = 154.694912
$ws.Cells.Item(4, 15).Value = 0.2098600871655206
$ws.Cells.Item(4, 16).Value = 0.2187357342528896
$ws.Cells.Item(4, 17).Value = 346.168694022912
$ws.Cells.Item(4, 18).Value = 3115.518246206208
$ws.Cells.Item(4, 19).Value = 0.06509524263421726
$ws.Cells.Item(4, 20).Value = 0.07068756633416415
$ws.Cells.Item(5, 7).Value = 6.713252999999999
$ws.Cells.Item(5, 8).Value = 20.139759
$ws.Cells.Item(5, 9).Value = 0.3101840064655811
$ws.Cells.Item(5, 10).Value = 0.3231642354899327
$ws.Cells.Item(5, 13).Value = 29.9106925
$ws.Cells.Item(5, 14).Value = 59.821385
$ws.Cells.Item(5, 15).Value = 0.1217310987299521
$ws.Cells.Item(5, 16).Value = 0.08458632803643724
$ws.Cells.Item(5, 17).Value = 200.7980461577025
$ws.Cells.Item(5, 18).Value = 1204.788276946215
$ws.Cells.Item(5, 19).Value = 0.03775903991551374
$ws.Cells.Item(5, 20).Value = 0.0273352760327959
$ws.Cells.Item(6, 7).Value = 6.713252999999999
$ws.Cells.Item(6, 8).Value = 20.139759
$ws.Cells.Item(6, 9).Value = 0.3101840064655811
$ws.Cells.Item(6, 10).Value = 0.3231642354899327
$ws.Cells.Item(6, 13).Value = 70.603826
$ws.Cells.Item(6, 14).Value = 211.811478
$ws.Cells.Item(6, 15).Value = 0.2873447785777061
$ws.Cells.Item(6, 16).Value = 0.2994974984278718
$ws.Cells.Item(6, 17).Value = 473.9813467059779
$ws.Cells.Item(6, 18).Value = 4265.832120353802
$ws.Cells.Item(6, 19).Value = 0.08912975465619818
$ws.Cells.Item(6, 20).Value = 0.09678688011059049
$ws.Cells.Item(7, 9).Value = 0.0154484264788496
$ws.Cells.Item(7, 10).Value = 0.01609489473505086
$ws.Cells.Item(7, 13).Value = 44.04223000000001
$ws.Cells.Item(7, 14).Value = 132.12669
$ws.Cells.Item(7, 15).Value = 0.1792438957545786
$ws.Cells.Item(7, 16).Value = 0.1868246872369915
$ws.Cells.Item(7, 17).Value = 14.72540215455334
$ws.Cells.Item(7, 18).Value = 132.52861939098
$ws.Cells.Item(7, 19).Value = 0.002769036145347189
$ws.Cells.Item(7, 20).Value = 0.003006923674988178
$ws.Cells.Item(8, 9).Value = 0.0154484264788496
$ws.Cells.Item(8, 10).Value = 0.01609489473505086
$ws.Cells.Item(8, 15).Value = 0.2018201397722426
$ws.Cells.Item(8, 16).Value = 0.2103557520458098
$ws.Cells.Item(8, 19).Value = 0.003117803591222639
$ws.Cells.Item(8, 20).Value = 0.003385653686089768
$ws.Cells.Item(9, 9).Value = 0.0154484264788496
$ws.Cells.Item(9, 10).Value = 0.01609489473505086
$ws.Cells.Item(9, 13).Value = 51.56497066666667
$ws.Cells.Item(9, 14).Value = 154.694912
$ws.Cells.Item(9, 15).Value = 0.2098600871655206
$ws.Cells.Item(9, 16).Value = 0.2187357342528896
$ws.Cells.Item(9, 17).Value = 17.24061043581156
$ws.Cells.Item(9, 18).Value = 155.165493922304
$ws.Cells.Item(9, 19).Value = 0.003242008127421512
$ws.Cells.Item(9, 20).Value = 0.003520528617594316
$ws.Cells.Item(10, 9).Value = 0.0154484264788496
$ws.Cells.Item(10, 10).Value = 0.01609489473505086
$ws.Cells.Item(10, 13).Value = 29.9106925
$ws.Cells.Item(10, 14).Value = 59.821385
$ws.Cells.Item(10, 15).Value = 0.1217310987299521
$ws.Cells.Item(10, 16).Value = 0.08458632803643724
$ws.Cells.Item(10, 17).Value = 10.00056027552833
$ws.Cells.Item(10, 18).Value = 60.00336165317
$ws.Cells.Item(10, 19).Value = 0.001880553928919246
$ws.Cells.Item(10, 20).Value = 0.001361408045770938
$ws.Cells.Item(11, 9).Value = 0.0154484264788496
$ws.Cells.Item(11, 10).Value = 0.01609489473505086
$ws.Cells.Item(11, 13).Value = 70.603826
$ws.Cells.Item(11, 14).Value = 211.811478
$ws.Cells.Item(11, 15).Value = 0.2873447785777061
$ws.Cells.Item(11, 16).Value = 0.2994974984278718
$ws.Cells.Item(11, 17).Value = 23.60620094623066
$ws.Cells.Item(11, 18).Value = 212.455808516076
$ws.Cells.Item(11, 19).Value = 0.004439024685939009
$ws.Cells.Item(11, 20).Value = 0.004820380710607655
$ws.Cells.Item(12, 7).Value = 6.661784666666667
$ws.Cells.Item(12, 8).Value = 19.985354
$ws.Cells.Item(12, 9).Value = 0.3078059262949933
$ws.Cells.Item(12, 10).Value = 0.3206866401135023
$ws.Cells.Item(12, 13).Value = 44.04223000000001
$ws.Cells.Item(12, 14).Value = 132.12669
$ws.Cells.Item(12, 15).Value = 0.1792438957545786
$ws.Cells.Item(12, 16).Value = 0.1868246872369915
$ws.Cells.Item(12, 17).Value = 293.3998524998067
$ws.Cells.Item(12, 18).Value = 2640.598672498261
$ws.Cells.Item(12, 19).Value = 0.0551723333654613
$ws.Cells.Item(12, 20).Value = 0.05991218124028672
$ws.Cells.Item(13, 7).Value = 6.661784666666667
$ws.Cells.Item(13, 8).Value = 19.985354
$ws.Cells.Item(13, 9).Value = 0.3078059262949933
$ws.Cells.Item(13, 10).Value = 0.3206866401135023
$ws.Cells.Item(13, 15).Value = 0.2018201397722426
$ws.Cells.Item(13, 16).Value = 0.2103557520458098
$ws.Cells.Item(13, 17).Value = 330.3543420063931
$ws.Cells.Item(13, 18).Value = 2973.189078057538
$ws.Cells.Item(13, 19).Value = 0.06212143506758014
$ws.Cells.Item(13, 20).Value = 0.06745827935211973
$ws.Cells.Item(14, 7).Value = 6.661784666666667
$ws.Cells.Item(14, 8).Value = 19.985354
$ws.Cells.Item(14, 9).Value = 0.3078059262949933
$ws.Cells.Item(14, 10).Value = 0.3206866401135023
$ws.Cells.Item(14, 13).Value = 51.56497066666667
$ws.Cells.Item(14, 14).Value = 154.694912
$ws.Cells.Item(14, 15).Value = 0.2098600871655206
$ws.Cells.Item(14, 16).Value = 0.2187357342528896
$ws.Cells.Item(14, 17).Value = 343.5147309243165
$ws.Cells.Item(14, 18).Value = 3091.632578318849
$ws.Cells.Item(14, 19).Value = 0.0645961785223311
$ws.Cells.Item(14, 20).Value = 0.07014562769031908
$ws.Cells.Item(15, 7).Value = 6.661784666666667
$ws.Cells.Item(15, 8).Value = 19.985354
$ws.Cells.Item(15, 9).Value = 0.3078059262949933
$ws.Cells.Item(15, 10).Value = 0.3206866401135023
$ws.Cells.Item(15, 13).Value = 29.9106925
$ws.Cells.Item(15, 14).Value = 59.821385
$ws.Cells.Item(15, 15).Value = 0.1217310987299521
$ws.Cells.Item(15, 16).Value = 0.08458632803643724
$ws.Cells.Item(15, 17).Value = 199.2585926658817
$ws.Cells.Item(15, 18).Value = 1195.55155599529
$ws.Cells.Item(15, 19).Value = 0.03746955360348018
$ws.Cells.Item(15, 20).Value = 0.0271257053375436
$ws.Cells.Item(16, 7).Value = 6.661784666666667
$ws.Cells.Item(16, 8).Value = 19.985354
$ws.Cells.Item(16, 9).Value = 0.3078059262949933
$ws.Cells.Item(16, 10).Value = 0.3206866401135023
$ws.Cells.Item(16, 13).Value = 70.603826
$ws.Cells.Item(16, 14).Value = 211.811478
$ws.Cells.Item(16, 15).Value = 0.2873447785777061
$ws.Cells.Item(16, 16).Value = 0.2994974984278718
$ws.Cells.Item(16, 17).Value = 470.3474854548013
$ws.Cells.Item(16, 18).Value = 4233.127369093212
$ws.Cells.Item(16, 19).Value = 0.08844642573614059
$ws.Cells.Item(16, 20).Value = 0.09604484649323312
$ws.Cells.Item(17, 7).Value = 2.607918
$ws.Cells.Item(17, 8).Value = 5.215835999999999
$ws.Cells.Item(17, 9).Value = 0.1204981331366039
$ws.Cells.Item(17, 10).Value = 0.08369373503331734
$ws.Cells.Item(17, 13).Value = 44.04223000000001
$ws.Cells.Item(17, 14).Value = 132.12669
$ws.Cells.Item(17, 15).Value = 0.1792438957545786
$ws.Cells.Item(17, 16).Value = 0.1868246872369915
$ws.Cells.Item(17, 17).Value = 114.85852437714
$ws.Cells.Item(17, 18).Value = 689.15114626284
$ws.Cells.Item(17, 19).Value = 0.02159855481455877
$ws.Cells.Item(17, 20).Value = 0.01563605587129515
$ws.Cells.Item(18, 7).Value = 2.607918
$ws.Cells.Item(18, 8).Value = 5.215835999999999
$ws.Cells.Item(18, 9).Value = 0.1204981331366039
$ws.Cells.Item(18, 10).Value = 0.08369373503331734
$ws.Cells.Item(18, 15).Value = 0.2018201397722426
$ws.Cells.Item(18, 16).Value = 0.2103557520458098
$ws.Cells.Item(18, 17).Value = 129.325260122482
$ws.Cells.Item(18, 18).Value = 775.9515607348919
$ws.Cells.Item(18, 19).Value = 0.02431895007192369
$ws.Cells.Item(18, 20).Value = 0.01760545857445621
$ws.Cells.Item(19, 7).Value = 2.607918
$ws.Cells.Item(19, 8).Value = 5.215835999999999
$ws.Cells.Item(19, 9).Value = 0.1204981331366039
$ws.Cells.Item(19, 10).Value = 0.08369373503331734
$ws.Cells.Item(19, 13).Value = 51.56497066666667
$ws.Cells.Item(19, 14).Value = 154.694912
$ws.Cells.Item(19, 15).Value = 0.2098600871655206
$ws.Cells.Item(19, 16).Value = 0.2187357342528896
$ws.Cells.Item(19, 17).Value = 134.477215171072
$ws.Cells.Item(19, 18).Value = 806.863291026432
$ws.Cells.Item(19, 19).Value = 0.02528774872333019
$ws.Cells.Item(19, 20).Value = 0.01830681058487946
$ws.Cells.Item(20, 7).Value = 2.607918
$ws.Cells.Item(20, 8).Value = 5.215835999999999
$ws.Cells.Item(20, 9).Value = 0.1204981331366039
$ws.Cells.Item(20, 10).Value = 0.08369373503331734
$ws.Cells.Item(20, 13).Value = 29.9106925
$ws.Cells.Item(20, 14).Value = 59.821385
$ws.Cells.Item(20, 15).Value = 0.1217310987299521
$ws.Cells.Item(20, 16).Value = 0.08458632803643724
$ws.Cells.Item(20, 17).Value = 78.00463336321499
$ws.Cells.Item(20, 18).Value = 312.01853345286
$ws.Cells.Item(20, 19).Value = 0.01466837014162684
$ws.Cells.Item(20, 20).Value = 0.007079345726122841
$ws.Cells.Item(21, 7).Value = 2.607918
$ws.Cells.Item(21, 8).Value = 5.215835999999999
$ws.Cells.Item(21, 9).Value = 0.1204981331366039
$ws.Cells.Item(21, 10).Value = 0.08369373503331734
$ws.Cells.Item(21, 13).Value = 70.603826
$ws.Cells.Item(21, 14).Value = 211.811478
$ws.Cells.Item(21, 15).Value = 0.2873447785777061
$ws.Cells.Item(21, 16).Value = 0.2994974984278718
$ws.Cells.Item(21, 17).Value = 184.128988694268
$ws.Cells.Item(21, 18).Value = 1104.773932165608
$ws.Cells.Item(21, 19).Value = 0.0346245093851644
$ws.Cells.Item(21, 20).Value = 0.02506606427656368
$ws.Cells.Item(22, 7).Value = 5.325505333333333
$ws.Cells.Item(22, 8).Value = 15.976516
$ws.Cells.Item(22, 9).Value = 0.2460635076239721
$ws.Cells.Item(22, 10).Value = 0.2563604946281968
$ws.Cells.Item(22, 13).Value = 44.04223000000001
$ws.Cells.Item(22, 14).Value = 132.12669
$ws.Cells.Item(22, 15).Value = 0.1792438957545786
$ws.Cells.Item(22, 16).Value = 0.1868246872369915
$ws.Cells.Item(22, 17).Value = 234.5471307568934
$ws.Cells.Item(22, 18).Value = 2110.924176812041
$ws.Cells.Item(22, 19).Value = 0.04410538170955722
$ws.Cells.Item(22, 20).Value = 0.04789446922883331
$ws.Cells.Item(23, 7).Value = 5.325505333333333
$ws.Cells.Item(23, 8).Value = 15.976516
$ws.Cells.Item(23, 9).Value = 0.2460635076239721
$ws.Cells.Item(23, 10).Value = 0.2563604946281968
$ws.Cells.Item(23, 15).Value = 0.2018201397722426
$ws.Cells.Item(23, 16).Value = 0.2103557520458098
$ws.Cells.Item(23, 17).Value = 264.0889638849835
$ws.Cells.Item(23, 18).Value = 2376.800674964852
$ws.Cells.Item(23, 19).Value = 0.04966057150151832
$ws.Cells.Item(23, 20).Value = 0.05392690464235012
$ws.Cells.Item(24, 7).Value = 5.325505333333333
$ws.Cells.Item(24, 8).Value = 15.976516
$ws.Cells.Item(24, 9).Value = 0.2460635076239721
$ws.Cells.Item(24, 10).Value = 0.2563604946281968
$ws.Cells.Item(24, 13).Value = 51.56497066666667
$ws.Cells.Item(24, 14).Value = 154.694912
$ws.Cells.Item(24, 15).Value = 0.2098600871655206
$ws.Cells.Item(24, 16).Value = 0.2187357342528896
$ws.Cells.Item(24, 17).Value = 274.6095262985103
$ws.Cells.Item(24, 18).Value = 2471.485736686592
$ws.Cells.Item(24, 19).Value = 0.05163890915822052
$ws.Cells.Item(24, 20).Value = 0.05607520102593257
$ws.Cells.Item(25, 7).Value = 5.325505333333333
$ws.Cells.Item(25, 8).Value = 15.976516
$ws.Cells.Item(25, 9).Value = 0.2460635076239721
$ws.Cells.Item(25, 10).Value = 0.2563604946281968
$ws.Cells.Item(25, 13).Value = 29.9106925
$ws.Cells.Item(25, 14).Value = 59.821385
$ws.Cells.Item(25, 15).Value = 0.1217310987299521
$ws.Cells.Item(25, 16).Value = 0.08458632803643724
$ws.Cells.Item(25, 17).Value = 159.2895524324433
$ws.Cells.Item(25, 18).Value = 955.73731459466
$ws.Cells.Item(25, 19).Value = 0.02995358114041206
$ws.Cells.Item(25, 20).Value = 0.02168459289420396
$ws.Cells.Item(26, 7).Value = 5.325505333333333
$ws.Cells.Item(26, 8).Value = 15.976516
$ws.Cells.Item(26, 9).Value = 0.2460635076239721
$ws.Cells.Item(26, 10).Value = 0.2563604946281968
$ws.Cells.Item(26, 13).Value = 70.603826
$ws.Cells.Item(26, 14).Value = 211.811478
$ws.Cells.Item(26, 15).Value = 0.2873447785777061
$ws.Cells.Item(26, 16).Value = 0.2994974984278718
$ws.Cells.Item(26, 17).Value = 376.0010519167386
$ws.Cells.Item(26, 18).Value = 3384.009467250648
$ws.Cells.Item(26, 19).Value = 0.07070506411426397
$ws.Cells.Item(26, 20).Value = 0.09678688011059049

Write-Output "Applied 278 cell updates"